$d = $word.ActiveDocument

# Locate the paragraph containing the "Ver no Jupiter..." text and the
# one right before it (the blank separator paragraph) plus the
# "© 2020 ..." footer paragraph that follows it, then remove all three
# as a single range so the "LOB1037: ..." paragraph is immediately
# followed by the trailing blank / page-break paragraphs again.
$count = $d.Paragraphs.Count
$startPara = $null
$endPara = $null

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        # the blank paragraph immediately preceding this one is part of
        # the block being removed
        $startPara = $i - 1
    }
    if ($t -like "*Powered by Jekyll*") {
        $endPara = $i
    }
}

$r = $d.Range($d.Paragraphs.Item($startPara).Range.Start, $d.Paragraphs.Item($endPara).Range.End)
$r.Delete()
